$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: was Java / 2545.0 / Lab 11 gui and database / Thu May 02 16:20:00 CDT 2019
# Becomes: CMST / 1234.0 / Essay / Fri May 10 00:00:00 CDT 2019
$ws.Range("A3").Value = "CMST"
$ws.Range("B3").Value = 1234.0
$ws.Range("C3").Value = "Essay"
$ws.Range("D3").Value = "Fri May 10 00:00:00 CDT 2019"

# Row 4: was CMST / 1234.0 / Essay / Fri May 10 00:00:00 CDT 2019
# Becomes: MATH / 3423.0 / Chapter 10 questions / Mon May 06 13:54:27 CDT 2019
$ws.Range("A4").Value = "MATH"
$ws.Range("B4").Value = 3423.0
$ws.Range("C4").Value = "Chapter 10 questions"
$ws.Range("D4").Value = "Mon May 06 13:54:27 CDT 2019"
